# Fix: convert string to double for price
# The values in column B (rows 2-7) were stored scaled by 1000; divide each
# by 1000 to get the correct price value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 616451
$ws.Range("B3").Value = 575514
$ws.Range("B4").Value = 8060647
$ws.Range("B5").Value = 469771
$ws.Range("B6").Value = 548604
$ws.Range("B7").Value = 3079553

# Move/restore the active selection to B7 (matches the saved view state).
$ws.Range("B7").Select()
